# Add a new snapshot column (K) to the LDLC price-tracking sheet.
# K1 gets a new timestamp header (styled like the other header cells),
# and K2:K100 get a copy of the latest prices currently in column J
# (rows 101+ have no price data yet, so column K stays blank there,
# exactly like column J does).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell K1 -------------------------------------------------
$j1 = $ws.Range("J1")
$k1 = $ws.Range("K1")

$k1.Value = "2026-01-28 00:51:24"

# Match the header formatting used by A1:J1 (bold font, thin box border,
# centered horizontally, top-aligned vertically).
$k1.Font.Bold = $j1.Font.Bold
$k1.HorizontalAlignment = $j1.HorizontalAlignment
$k1.VerticalAlignment = $j1.VerticalAlignment
$k1.Borders.LineStyle = $j1.Borders.LineStyle

# --- Data rows: copy column J's current prices into column K --------
$lastRow = 204
for ($r = 2; $r -le $lastRow; $r++) {
    $src = $ws.Cells.Item($r, 10)   # column J
    $val = $src.Value2
    if ($val -ne $null -and $val -ne "") {
        $ws.Cells.Item($r, 11).Value = $val   # column K
    }
}
